$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sending cluster / Ligand symbol / Receptor symbol / Target cluster grid for Vegfc-Kdr
# now spans ECs/FAPs/sCs x ECs/FAPs/sCs (3x3), per updated NATMI run (Dr Hou advice).
$data = @(
    @("ECs", "Vegfc", "Kdr", "ECs", 2, 0.6666666666666666, 1.164278333333333, 3.492835, 0.2070506538112546, 0.2070506538112546, 3, 1, 62.351132, 187.053396, 0.9620005273240912, 0.9620005273240912, 72.59407204640667, 653.3466484176599, 0.1991828381492248, 0.1991828381492248),
    @("ECs", "Vegfc", "Kdr", "FAPs", 2, 0.6666666666666666, 1.164278333333333, 3.492835, 0.2070506538112546, 0.2070506538112546, 3, 1, 0.8562703333333334, 2.568811, 0.01321118776478095, 0.01321118776478095, 0.9969369965761111, 8.972432969185, 0.002735385064321144, 0.002735385064321144),
    @("ECs", "Vegfc", "Kdr", "sCs", 2, 0.6666666666666666, 1.164278333333333, 3.492835, 0.2070506538112546, 0.2070506538112546, 3, 1, 1.606628666666667, 4.819886, 0.02478828491112776, 0.02478828491112776, 1.870562946312222, 16.83506651681, 0.005132430597708661, 0.005132430597708661),
    @("FAPs", "Vegfc", "Kdr", "ECs", 3, 1, 3.886460333333333, 11.659381, 0.6911527338349851, 0.6911527338349851, 3, 1, 62.351132, 187.053396, 0.9620005273240912, 0.9620005273240912, 242.3252012564307, 2180.926811307876, 0.664889294410743, 0.664889294410743),
    @("FAPs", "Vegfc", "Kdr", "FAPs", 3, 1, 3.886460333333333, 11.659381, 0.6911527338349851, 0.6911527338349851, 3, 1, 0.8562703333333334, 2.568811, 0.01321118776478095, 0.01321118776478095, 3.327860685110111, 29.950746165991, 0.009130948540835663, 0.009130948540835663),
    @("FAPs", "Vegfc", "Kdr", "sCs", 3, 1, 3.886460333333333, 11.659381, 0.6911527338349851, 0.6911527338349851, 3, 1, 1.606628666666667, 4.819886, 0.02478828491112776, 0.02478828491112776, 6.244098583396222, 56.196887250566, 0.01713249088340646, 0.01713249088340646),
    @("sCs", "Vegfc", "Kdr", "ECs", 3, 1, 0.5724183333333334, 1.717255, 0.1017966123537602, 0.1017966123537602, 3, 1, 62.351132, 187.053396, 0.9620005273240912, 0.9620005273240912, 35.69093106088667, 321.21837954798, 0.09792839476412345, 0.09792839476412343),
    @("sCs", "Vegfc", "Kdr", "FAPs", 3, 1, 0.5724183333333334, 1.717255, 0.1017966123537602, 0.1017966123537602, 3, 1, 0.8562703333333334, 2.568811, 0.01321118776478095, 0.01321118776478095, 0.4901448370894445, 4.411303533805, 0.001344854159624147, 0.001344854159624147),
    @("sCs", "Vegfc", "Kdr", "sCs", 3, 1, 0.5724183333333334, 1.717255, 0.1017966123537602, 0.1017966123537602, 3, 1, 1.606628666666667, 4.819886, 0.02478828491112776, 0.02478828491112776, 0.9196637036588889, 8.27697333293, 0.002523363430012636, 0.002523363430012636)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowVals = $data[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowVals[$c]
    }
}
